$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised regression-table values for "Model 1c" (R&R Democratization).
# Row 2 ((Intercept)) and row 3 (lrscale) are updated in place; rows 4-10
# are newly added variables. Column C ("icc") holds a numeric-looking
# string value, so it is entered with a leading apostrophe to keep it
# stored as text (matching the workbook's existing inlineStr/text usage).

$rows = @(
    @{ r = 2;  A = "(Intercept)";   B = 1.912;  C = "0.0785632550862939"; D = 0.07691966893576119;  E = 32130; F = "Model 1c"; G = 1.785467144600673;   H = 2.038532855399327;   I = 12 },
    @{ r = 3;  A = "lrscale";       B = -0.325; C = "0.0785632550862939"; D = 0.02219912530737515;  E = 32130; F = "Model 1c"; G = -0.3615175611306322; H = -0.2884824388693679; I = 12 },
    @{ r = 4;  A = "age";           B = 0.065;  C = "0.0785632550862939"; D = 0.02741684820586976;  E = 32130; F = "Model 1c"; G = 0.01989928470134425;  H = 0.1101007152986558;  I = 12 },
    @{ r = 5;  A = "educ";          B = 0.387;  C = "0.0785632550862939"; D = 0.01456523292771531;  E = 32130; F = "Model 1c"; G = 0.3630401918339083;   H = 0.4109598081660917;  I = 12 },
    @{ r = 6;  A = "polint";        B = -0.789; C = "0.0785632550862939"; D = 0.01714843936538139;  E = 32130; F = "Model 1c"; G = -0.8172091827560524;  H = -0.7607908172439477; I = 12 },
    @{ r = 7;  A = "sexMale";       B = 0.05;   C = "0.0785632550862939"; D = 0.009844321487365165; E = 32130; F = "Model 1c"; G = 0.03380609115328431;  H = 0.06619390884671569; I = 12 },
    @{ r = 8;  A = "surveyevs2008"; B = -0.139; C = "0.0785632550862939"; D = 0.01264080065724252;  E = 32130; F = "Model 1c"; G = -0.1597941170811639;  H = -0.1182058829188361; I = 12 },
    @{ r = 9;  A = "surveywvs1994"; B = -0.247; C = "0.0785632550862939"; D = 0.02088511705940227;  E = 32130; F = "Model 1c"; G = -0.2813560175627167;  H = -0.2126439824372833; I = 12 },
    @{ r = 10; A = "surveywvs2005"; B = 0.304;  C = "0.0785632550862939"; D = 0.01654853406312838;  E = 32130; F = "Model 1c"; G = 0.2767776614661538;   H = 0.3312223385338462;  I = 12 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = "'" + $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
}
